$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.436.04'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '1.823.81'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.94'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5117'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -4.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3916'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07645'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.61'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.99'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.267'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.003'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.506'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '1.821.91'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.32'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001094'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06682'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.67'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.152'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("D23").Value = '28.444.85'
$ws.Range("E23").Value = '  -0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.17'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("E25").Value = '  +7.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.77'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.91'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = '2.031.47'
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.379'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.107'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1086'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.642'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.664'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07060'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2199'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02317'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.824'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.164'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6245'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.22'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.172'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.36'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5878'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.13'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.979'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.197'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06908'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.30%  '
